# "update template fill kabupaten"
#
# The sheet holds a lookup table: column A = kabupaten/kota name, column B =
# id_provinsi, and (previously) column C held a helper formula
#   =PROPER(SUBSTITUTE(SUBSTITUTE(UPPER(A#), "KAB.", "Kabupaten"), "KAB ", "Kabupaten "))
# that was used while drafting the list but was never filled in for rows
# 493-570 (column A/B were empty there).
#
# This edit fills in the missing kabupaten/kota entries for Aceh province
# (id_provinsi = 11) directly into columns A and B for rows 493-515, and
# removes the now-unused helper formula from column C entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "Kabupaten Aceh Selatan",
    "Kabupaten Aceh Tenggara",
    "Kabupaten Aceh Timur",
    "Kabupaten Aceh Tengah",
    "Kabupaten Aceh Barat",
    "Kabupaten Aceh Besar",
    "Kabupaten Pidie",
    "Kabupaten Aceh Utara",
    "Kabupaten Simeulue",
    "Kabupaten Aceh Singkil",
    "Kabupaten Bireuen",
    "Kabupaten Aceh Barat Daya",
    "Kabupaten Gayo Lues",
    "Kabupaten Aceh Jaya",
    "Kabupaten Nagan Raya",
    "Kabupaten Aceh Tamiang",
    "Kabupaten Bener Meriah",
    "Kabupaten Pidie Jaya",
    "Kota Banda Aceh",
    "Kota Sabang",
    "Kota Lhokseumawe",
    "Kota Langsa",
    "Kota Subulussalam"
)

$startRow = 493
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = 11
}

# Drop the helper PROPER/SUBSTITUTE formulas that used to live in column C
# for the whole table (rows 1-570) - they're no longer needed now that the
# names are typed straight into column A.
$ws.Range("C1:C570").ClearContents()

# Leave the selection where the editor ended up.
$ws.Range("A502").Select()
